# Revert "update Todo/Bug form"
# - Drop the "Bug List" sheet entirely
# - Rename "Todo List" back to "Sheet1"
# - Remove the "Status" column data (B1, B3:B10) added to the Todo sheet,
#   along with the special header formatting (fill/border/bold/height) that
#   was introduced for that header row
# - Undo the frozen header pane / selection on the Todo sheet
# - Drop the stray "D31" selection that had been left on the old "Sheet2"

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$todo = $wb.Worksheets.Item("Todo List")

# Restore the original (unstyled) look of the header cells, and drop the
# custom row height / thick bottom border that came with that styling.
$todo.Range("A1:B2").Style = "Normal"
$todo.Rows.Item(1).AutoFit()
$todo.Rows.Item(2).AutoFit()

# Remove the "Status" column values that were introduced alongside the Bug
# List sheet; keep the B2 "Status" header text itself.
$todo.Range("B1").ClearContents()
$todo.Range("B3:B10").ClearContents()

# Undo the frozen pane / split selection and go back to a plain selection.
$todo.Activate()
$excel.ActiveWindow.FreezePanes = $false
$todo.Range("A10").Select()

# The sheet was originally named "Sheet1".
$todo.Name = "Sheet1"

# Drop the "Bug List" sheet that this commit introduced.
$wb.Worksheets.Item("Bug List").Delete()

# The old "Sheet2" tab had picked up a stray selection at D31; clear it back
# to the default top-left selection.
$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet2.Activate()
$sheet2.Range("A1").Select()

$wb.Worksheets.Item("Sheet1").Activate()
